$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'67.399.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "'  +0.80%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value2 = "'3.834.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "'  +0.89%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value2 = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "'  -0.21%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value2 = "'456.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "'  +8.31%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value2 = "'146.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "'  +14.33%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value2 = "'  +3.49%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value2 = "'  -0.11%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value2 = "'0.741"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "'  +3.37%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value2 = "'  -3.38%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value2 = "'0.0000318"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "'  -8.92%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value2 = "'43.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "'  +8.82%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value2 = "'10.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "'  +2.55%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value2 = "'4.442.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "'  +0.64%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value2 = "'14.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "'  -4.83%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value2 = "'WrappedEther"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value2 = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value2 = "'3.832.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "'  +1.14%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value2 = "'TRON"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value2 = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value2 = "'0.137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "'  -0.31%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value2 = "'20.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "'  +3.27%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value2 = "'1.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "'  +8.62%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value2 = "'67.433.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "'  +0.62%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value2 = "'433.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "'  +6.92%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value2 = "'14.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "'  +0.27%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value2 = "'  +8.74%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value2 = "'86.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "'  +3.77%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value2 = "'  +9.71%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value2 = "'10.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "'  +14.86%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value2 = "'37.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "'  +0.79%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value2 = "'5.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "'  +0.09%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value2 = "'9.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "'  +2.47%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value2 = "'733.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "'  +1.17%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value2 = "'13.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "'  +11.53%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value2 = "'  +12.10%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value2 = "'2.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "'  -1.68%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value2 = "'43.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "'  +12.00%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value2 = "'0.161"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "'  +5.80%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value2 = "'57.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "'  +4.37%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value2 = "'Dai"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value2 = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value2 = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "'  +0.15%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value2 = "'NEARProtocol"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value2 = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value2 = "'5.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "'  +3.48%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value2 = "'0.0475"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "'  +5.87%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value2 = "'0.352"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "'  +12.73%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value2 = "'  -0.42%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value2 = "'PEPE"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value2 = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value2 = "'0.0₃0690"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "'  -7.51%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value2 = "'Fetch.AI"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value2 = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value2 = "'2.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "'  +15.40%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value2 = "'  -0.28%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value2 = "'  +5.21%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value2 = "'3.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "'  +4.20%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value2 = "'  +5.62%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value2 = "'2.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "'  +5.41%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value2 = "'2.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "'  +5.42%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value2 = "'143.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "'  +0.20%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value2 = "'  +3.01%  "
$ws.Range("E51").Style = "Normal"
